$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INDUSIND")

# Row 7
$ws.Range("F7").Value = 1538.6
$ws.Range("G7").Value = 1574
$ws.Range("H7").Value = 1527.55
$ws.Range("I7").Value = 1570
$ws.Range("J7").Value = 1556.15

# Row 9
$ws.Range("G9").Value = 1555.5
$ws.Range("H9").Value = 1533.55
$ws.Range("I9").Value = 1534.6

# Row 10
$ws.Range("G10").Value = 1541.5
$ws.Range("H10").Value = 1529.05
$ws.Range("I10").Value = 1531.35

# Row 11
$ws.Range("G11").Value = 1541.35
$ws.Range("H11").Value = 1527.55
$ws.Range("I11").Value = 1532.5

# Row 12
$ws.Range("G12").Value = 1537
$ws.Range("H12").Value = 1531
$ws.Range("I12").Value = 1535.15

# Row 13
$ws.Range("G13").Value = 1543.1
$ws.Range("H13").Value = 1534
$ws.Range("I13").Value = 1539.4

# Row 14
$ws.Range("G14").Value = 1542.4
$ws.Range("H14").Value = 1536.45
$ws.Range("I14").Value = 1539.1

# Row 15
$ws.Range("G15").Value = 1547.5
$ws.Range("H15").Value = 1539.1
$ws.Range("I15").Value = 1546.35

# Row 16
$ws.Range("G16").Value = 1552
$ws.Range("H16").Value = 1544.1
$ws.Range("I16").Value = 1545.85

# Row 17
$ws.Range("G17").Value = 1551.95
$ws.Range("H17").Value = 1545.1
$ws.Range("I17").Value = 1551.65

# Row 18
$ws.Range("G18").Value = 1552.9
$ws.Range("H18").Value = 1546.1
$ws.Range("I18").Value = 1548.25

# Row 19
$ws.Range("G19").Value = 1551.3
$ws.Range("H19").Value = 1540.4
$ws.Range("I19").Value = 1550.5

# Row 20
$ws.Range("G20").Value = 1573.95
$ws.Range("H20").Value = 1549.85
$ws.Range("I20").Value = 1564.55

# Row 21
$ws.Range("G21").Value = 1574
$ws.Range("H21").Value = 1563.05
$ws.Range("I21").Value = 1570
